$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Egf"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1243246666666667
$ws.Range("H2").Value = 0.372974
$ws.Range("I2").Value = 0.09963085929726231
$ws.Range("J2").Value = 0.09963085929726233
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1680563333333333
$ws.Range("N2").Value = 0.504169
$ws.Range("O2").Value = 0.03491515448966758
$ws.Range("P2").Value = 0.03491515448966757
$ws.Range("Q2").Value = 0.02089354762288889
$ws.Range("R2").Value = 0.188041928606
$ws.Range("S2").Value = 0.003478626844302247
$ws.Range("T2").Value = 0.003478626844302247

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Egf"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1243246666666667
$ws.Range("H3").Value = 0.372974
$ws.Range("I3").Value = 0.09963085929726231
$ws.Range("J3").Value = 0.09963085929726233
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3027816666666667
$ws.Range("N3").Value = 0.908345
$ws.Range("O3").Value = 0.06290550590162643
$ws.Range("P3").Value = 0.06290550590162643
$ws.Range("Q3").Value = 0.03764322978111111
$ws.Range("R3").Value = 0.33878906803
$ws.Range("S3").Value = 0.006267329607508047
$ws.Range("T3").Value = 0.006267329607508048

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Egf"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1243246666666667
$ws.Range("H4").Value = 0.372974
$ws.Range("I4").Value = 0.09963085929726231
$ws.Range("J4").Value = 0.09963085929726233
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3920839999999999
$ws.Range("N4").Value = 1.176252
$ws.Range("O4").Value = 0.08145883681618755
$ws.Range("P4").Value = 0.08145883681618755
$ws.Range("Q4").Value = 0.04874571260533333
$ws.Range("R4").Value = 0.438711413448
$ws.Range("S4").Value = 0.008115813909352234
$ws.Range("T4").Value = 0.008115813909352235

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Egf"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1243246666666667
$ws.Range("H5").Value = 0.372974
$ws.Range("I5").Value = 0.09963085929726231
$ws.Range("J5").Value = 0.09963085929726233
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.950355666666666
$ws.Range("N5").Value = 11.851067
$ws.Range("O5").Value = 0.8207205027925185
$ws.Range("P5").Value = 0.8207205027925184
$ws.Range("Q5").Value = 0.4911266514731111
$ws.Range("R5").Value = 4.420139863258
$ws.Range("S5").Value = 0.0817690889360998
$ws.Range("T5").Value = 0.0817690889360998

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Egf"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7328223333333334
$ws.Range("H6").Value = 2.198467
$ws.Range("I6").Value = 0.5872665557027417
$ws.Range("J6").Value = 0.5872665557027417
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1680563333333333
$ws.Range("N6").Value = 0.504169
$ws.Range("O6").Value = 0.03491515448966758
$ws.Range("P6").Value = 0.03491515448966757
$ws.Range("Q6").Value = 0.1231554343247778
$ws.Range("R6").Value = 1.108398908923
$ws.Range("S6").Value = 0.0205045025189762
$ws.Range("T6").Value = 0.02050450251897619

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Egf"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7328223333333334
$ws.Range("H7").Value = 2.198467
$ws.Range("I7").Value = 0.5872665557027417
$ws.Range("J7").Value = 0.5872665557027417
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3027816666666667
$ws.Range("N7").Value = 0.908345
$ws.Range("O7").Value = 0.06290550590162643
$ws.Range("P7").Value = 0.06290550590162643
$ws.Range("Q7").Value = 0.2218851674572222
$ws.Range("R7").Value = 1.996966507115
$ws.Range("S7").Value = 0.03694229978558664
$ws.Range("T7").Value = 0.03694229978558664

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Egf"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7328223333333334
$ws.Range("H8").Value = 2.198467
$ws.Range("I8").Value = 0.5872665557027417
$ws.Range("J8").Value = 0.5872665557027417
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3920839999999999
$ws.Range("N8").Value = 1.176252
$ws.Range("O8").Value = 0.08145883681618755
$ws.Range("P8").Value = 0.08145883681618755
$ws.Range("Q8").Value = 0.2873279117426666
$ws.Range("R8").Value = 2.585951205684
$ws.Range("S8").Value = 0.04783805052859415
$ws.Range("T8").Value = 0.04783805052859415

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Egf"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7328223333333334
$ws.Range("H9").Value = 2.198467
$ws.Range("I9").Value = 0.5872665557027417
$ws.Range("J9").Value = 0.5872665557027417
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.950355666666666
$ws.Range("N9").Value = 11.851067
$ws.Range("O9").Value = 0.8207205027925185
$ws.Range("P9").Value = 0.8207205027925184
$ws.Range("Q9").Value = 2.894908857143222
$ws.Range("R9").Value = 26.054179714289
$ws.Range("S9").Value = 0.4819817028695847
$ws.Range("T9").Value = 0.4819817028695846

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Egf"
$ws.Range("C10").Value = "Erbb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2213483333333333
$ws.Range("H10").Value = 0.664045
$ws.Range("I10").Value = 0.1773833402919521
$ws.Range("J10").Value = 0.1773833402919521
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1680563333333333
$ws.Range("N10").Value = 0.504169
$ws.Range("O10").Value = 0.03491515448966758
$ws.Range("P10").Value = 0.03491515448966757
$ws.Range("Q10").Value = 0.03719898928944445
$ws.Range("R10").Value = 0.334790903605
$ws.Range("S10").Value = 0.006193366730186785
$ws.Range("T10").Value = 0.006193366730186784

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Egf"
$ws.Range("C11").Value = "Erbb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2213483333333333
$ws.Range("H11").Value = 0.664045
$ws.Range("I11").Value = 0.1773833402919521
$ws.Range("J11").Value = 0.1773833402919521
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3027816666666667
$ws.Range("N11").Value = 0.908345
$ws.Range("O11").Value = 0.06290550590162643
$ws.Range("P11").Value = 0.06290550590162643
$ws.Range("Q11").Value = 0.06702021728055556
$ws.Range("R11").Value = 0.603181955525
$ws.Range("S11").Value = 0.0111583887595856
$ws.Range("T11").Value = 0.0111583887595856

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Egf"
$ws.Range("C12").Value = "Erbb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2213483333333333
$ws.Range("H12").Value = 0.664045
$ws.Range("I12").Value = 0.1773833402919521
$ws.Range("J12").Value = 0.1773833402919521
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.3920839999999999
$ws.Range("N12").Value = 1.176252
$ws.Range("O12").Value = 0.08145883681618755
$ws.Range("P12").Value = 0.08145883681618755
$ws.Range("Q12").Value = 0.08678713992666666
$ws.Range("R12").Value = 0.7810842593399999
$ws.Range("S12").Value = 0.0144494405707524
$ws.Range("T12").Value = 0.0144494405707524

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Egf"
$ws.Range("C13").Value = "Erbb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2213483333333333
$ws.Range("H13").Value = 0.664045
$ws.Range("I13").Value = 0.1773833402919521
$ws.Range("J13").Value = 0.1773833402919521
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.950355666666666
$ws.Range("N13").Value = 11.851067
$ws.Range("O13").Value = 0.8207205027925185
$ws.Range("P13").Value = 0.8207205027925184
$ws.Range("Q13").Value = 0.8744046428905555
$ws.Range("R13").Value = 7.869641786014999
$ws.Range("S13").Value = 0.1455821442314273
$ws.Range("T13").Value = 0.1455821442314273

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Egf"
$ws.Range("C14").Value = "Erbb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1693576666666667
$ws.Range("H14").Value = 0.508073
$ws.Range("I14").Value = 0.1357192447080439
$ws.Range("J14").Value = 0.1357192447080439
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1680563333333333
$ws.Range("N14").Value = 0.504169
$ws.Range("O14").Value = 0.03491515448966758
$ws.Range("P14").Value = 0.03491515448966757
$ws.Range("Q14").Value = 0.02846162848188889
$ws.Range("R14").Value = 0.256154656337
$ws.Range("S14").Value = 0.004738658396202351
$ws.Range("T14").Value = 0.00473865839620235

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Egf"
$ws.Range("C15").Value = "Erbb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1693576666666667
$ws.Range("H15").Value = 0.508073
$ws.Range("I15").Value = 0.1357192447080439
$ws.Range("J15").Value = 0.1357192447080439
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.3027816666666667
$ws.Range("N15").Value = 0.908345
$ws.Range("O15").Value = 0.06290550590162643
$ws.Range("P15").Value = 0.06290550590162643
$ws.Range("Q15").Value = 0.05127839657611111
$ws.Range("R15").Value = 0.461505569185
$ws.Range("S15").Value = 0.008537487748946134
$ws.Range("T15").Value = 0.008537487748946134

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Egf"
$ws.Range("C16").Value = "Erbb3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1693576666666667
$ws.Range("H16").Value = 0.508073
$ws.Range("I16").Value = 0.1357192447080439
$ws.Range("J16").Value = 0.1357192447080439
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3920839999999999
$ws.Range("N16").Value = 1.176252
$ws.Range("O16").Value = 0.08145883681618755
$ws.Range("P16").Value = 0.08145883681618755
$ws.Range("Q16").Value = 0.06640243137733332
$ws.Range("R16").Value = 0.5976218823959999
$ws.Range("S16").Value = 0.01105553180748877
$ws.Range("T16").Value = 0.01105553180748877

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Egf"
$ws.Range("C17").Value = "Erbb3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1693576666666667
$ws.Range("H17").Value = 0.508073
$ws.Range("I17").Value = 0.1357192447080439
$ws.Range("J17").Value = 0.1357192447080439
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.950355666666666
$ws.Range("N17").Value = 11.851067
$ws.Range("O17").Value = 0.8207205027925185
$ws.Range("P17").Value = 0.8207205027925184
$ws.Range("Q17").Value = 0.669023018210111
$ws.Range("R17").Value = 6.021207163890999
$ws.Range("S17").Value = 0.1113875667554066
$ws.Range("T17").Value = 0.1113875667554066
